# LOT2043.xlsx update
# - Fill in the previously-missing course-plan data (objectives, teacher
#   names, syllabi, method, criteria, recovery rule, bibliography).
# - Two new rows are needed to hold both teacher names ("Docentes
#   responsáveis:" used to have only one name where two belong), so rows
#   13-14 are inserted and the remainder of the table shifts down by two
#   (old row 13 -> 15, ... old row 21 -> 23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two "Docentes responsáveis" rows; everything at/after
# row 13 (and its row height/formatting) shifts down to row 15 onward.
$ws.Rows("13:14").Insert()

# The row-insert copies column A's formatting down into A13/A14 even
# though those cells should stay empty (no "label" in column A for these
# two rows) - clear that out.
$ws.Range("A13:A14").Clear()

# ...and it does not give B13:C14 the normal "data column" formatting
# (style 2 = black/wrap for column B, style 3 = red/wrap for column C)
# that every other data row uses, so copy it down from a row that already
# has it before filling in the values.
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C14").PasteSpecial(-4122)  # xlPasteFormats

# Row 10 (Objetivos:) previously had a teacher name left in column B/C by
# mistake; replace with the real Portuguese objectives text.
$ws.Range("B10:C10").Value = "Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia."

# Rows 13-14 (new, blank): the two responsible teachers.
$ws.Range("B13:C13").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("B14:C14").Value = "5817181 - Valdeir Arantes"

# Row 15 (Programa resumido:) gets the actual short syllabus text.
$ws.Range("B15:C15").Value = "1. Histórico da Engenharia Bioquímica 2. Engenharia Bioquímica: Definições e conceitos 3. Mercado de trabalho de Engenharia4. Áreas de atuação do Engenheiro Bioquímico5. A Indústria de Bioprocessos 6. Escalas de produção 7. Estudo de casos (processos biotecnológicos)8. Visita supervisionada."

# Row 17 (Programa:) gets the full syllabus text.
$ws.Range("B17:C17").Value = "1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso."

# Row 20 (Método:) gets the real teaching-method description.
$ws.Range("B20:C20").Value = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico."

# Row 21 (Critério:) gets the grading-criteria description.
$ws.Range("B21:C21").Value = "A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros."

# Row 22 (Norma de recuperação:) gets the recovery-grade formula.
$ws.Range("B22:C22").Value = "Média Final = (N + Prova Recuperação)/2"

# Row 23 (Bibliografia:) gets the reading list.
$ws.Range("B23:C23").Value = "Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial – EngenhariaBioquímica, vol. 2, São Paulo: Edgard Blücher, 2001.Shuler, L. M.; Kargi, F. Bioprocess Engineering – Basic Concepts. Second edition. NewJersey: PrenticeHall,2002.Arigos atuais relacionaos com o tema de Engenharia Bioquímica"
